# Applies the multiplication-table answer replacements described by the diff.
$d = $word.ActiveDocument

$replacements = @(
    @("93×75=6975", "12×60=720"),
    @("12×39=468", "60×20=1200"),
    @("72×68=4896", "22×81=1782"),
    @("87×15=1305", "54×87=4698"),
    @("20×96=1920", "60×94=5640"),
    @("54×12=648", "61×15=915"),
    @("78×96=7488", "68×16=1088"),
    @("21×97=2037", "29×11=319"),
    @("75×64=4800", "85×88=7480"),
    @("52×86=4472", "82×89=7298"),
    @("43×87=3741", "88×63=5544"),
    @("44×67=2948", "49×77=3773"),
    @("21×96=2016", "21×66=1386"),
    @("54×41=2214", "15×67=1005"),
    @("62×98=6076", "62×36=2232"),
    @("15×47=705", "48×99=4752"),
    @("81×95=7695", "12×75=900"),
    @("30×27=810", "19×79=1501"),
    @("66×66=4356", "86×32=2752"),
    @("53×83=4399", "61×73=4453"),
    @("13×82=1066", "89×85=7565"),
    @("66×99=6534", "82×53=4346"),
    @("88×66=5808", "73×85=6205"),
    @("99×17=1683", "27×22=594"),
    @("42×27=1134", "29×61=1769")
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $range = $d.Content
    $range.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)
}

$d.Save()
